$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.782.51"
$ws.Range("E2").Value = "  +4.78%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.767.52"
$ws.Range("E3").Value = "  +4.55%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.65"
$ws.Range("E5").Value = "  +0.98%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.00"
$ws.Range("E6").Value = "  +9.61%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.22%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.611"
$ws.Range("E8").Value = "  +2.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.763.21"
$ws.Range("E9").Value = "  +3.74%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.75"
$ws.Range("E10").Value = "  +2.66%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.113"
$ws.Range("E11").Value = "  +6.21%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.394"
$ws.Range("E12").Value = "  +3.69%  "

$ws.Range("E13").Value = "  +1.78%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.248.13"
$ws.Range("E14").Value = "  +3.74%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.31"
$ws.Range("E15").Value = "  +4.42%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.773.39"
$ws.Range("E16").Value = "  +4.79%  "

$ws.Range("E17").Value = "  +8.37%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.762.20"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.14"
$ws.Range("E19").Value = "  +4.98%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.95"
$ws.Range("E20").Value = "  +4.64%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "363.19"
$ws.Range("E21").Value = "  +3.54%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.04"
$ws.Range("E22").Value = "  +2.56%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("E24").Value = "  +1.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.07"
$ws.Range("E25").Value = "  +4.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.171"
$ws.Range("E26").Value = "  +5.92%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.62"
$ws.Range("E27").Value = "  +4.75%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.996"
$ws.Range("E28").Value = "  -0.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0916"
$ws.Range("E29").Value = "  +13.60%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.02"
$ws.Range("E30").Value = "  +1.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.20"
$ws.Range("E31").Value = "  +5.66%  "

$ws.Range("E32").Value = "  +19.77%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "174.83"
$ws.Range("E33").Value = "  +7.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.996"
$ws.Range("E34").Value = "  -0.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "20.68"
$ws.Range("E35").Value = "  +3.84%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.91"
$ws.Range("E36").Value = "  +5.74%  "

$ws.Range("E37").Value = "  +9.62%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.82"
$ws.Range("E38").Value = "  +9.64%  "

$ws.Range("E39").Value = "  +11.86%  "

$ws.Range("E40").Value = "  +5.34%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "341.60"
$ws.Range("E41").Value = "  +0.41%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.40"
$ws.Range("E42").Value = "  +2.36%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.88"
$ws.Range("E43").Value = "  +11.98%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.04"
$ws.Range("E44").Value = "  +8.57%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.24"
$ws.Range("E45").Value = "  +8.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0601"
$ws.Range("E46").Value = "  +6.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.653"
$ws.Range("E47").Value = "  +4.65%  "

$ws.Range("E48").Value = "  +4.64%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "138.67"
$ws.Range("E49").Value = "  +4.19%  "

$ws.Range("E50").Value = "  +2.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.997"
$ws.Range("E51").Value = "  -0.07%  "
